$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 257.42856
$ws.Range("J9").Value = 190.2
$ws.Range("L9").Value = 190.2
$ws.Range("N9").Value = -528.2

$ws.Range("H12").Value = 2358.2856
$ws.Range("I12").Value = 220
$ws.Range("J12").Value = 3213.6
$ws.Range("K12").Value = 220
$ws.Range("L12").Value = 3213.6
$ws.Range("M12").Value = -50
$ws.Range("N12").Value = -3553.6

$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("L13").ClearContents()

$ws.Range("H33").Value = 374.0811
$ws.Range("I33").Value = 195.12122
$ws.Range("K33").Value = 195.12122
$ws.Range("M33").Value = 33.87878000000001

$ws.Range("H40").Value = 2493.875
$ws.Range("I40").Value = 1750.5
$ws.Range("K40").Value = 1750.5
$ws.Range("M40").Value = -1575.5

$ws.Range("H42").Value = 113.6
$ws.Range("I42").Value = 56
$ws.Range("K42").Value = 168
$ws.Range("M42").Value = 62

$ws.Range("H48").Value = 2000
$ws.Range("I48").Value = 1000
$ws.Range("J48").Value = 3000
$ws.Range("K48").Value = 3000
$ws.Range("L48").Value = 9000
$ws.Range("M48").Value = -2708
$ws.Range("N48").Value = -9584

$ws.Range("H56").Value = 2000
$ws.Range("I56").Value = 1000
$ws.Range("J56").Value = 3000
$ws.Range("K56").Value = 3000
$ws.Range("L56").Value = 9000
$ws.Range("M56").Value = -2466
$ws.Range("N56").Value = -10068

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1956.75
$ws.Range("I61").Value = 1890
$ws.Range("J61").Value = 2157
$ws.Range("K61").Value = 1890
$ws.Range("L61").Value = 2157
$ws.Range("M61").Value = -1678
$ws.Range("N61").Value = -2581

$ws.Range("H132").Value = 2873.238
$ws.Range("I132").Value = 1904
$ws.Range("K132").Value = 5712
$ws.Range("M132").Value = -3182

$ws.Range("H136").Value = 1956.75
$ws.Range("I136").Value = 1890
$ws.Range("J136").Value = 2157
$ws.Range("K136").Value = 5670
$ws.Range("L136").Value = 6471
$ws.Range("M136").Value = -3120
$ws.Range("N136").Value = -11571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 638.8
$ws.Range("I22").Value = 638.8
$ws.Range("K22").Value = 638.8
$ws.Range("M22").Value = -465.8

$ws.Range("H134").Value = 854782.75
$ws.Range("I134").Value = 1179915.9
$ws.Range("J134").Value = 4434.4614
$ws.Range("K134").Value = 3539747.7
$ws.Range("L134").Value = 13303.3842
$ws.Range("M134").Value = -3537212.7
$ws.Range("N134").Value = -18373.3842

$ws.Range("H140").Value = 19452.5
$ws.Range("J140").Value = 19452.5
$ws.Range("L140").Value = 19452.5
$ws.Range("N140").Value = -29812.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 26822.4
$ws.Range("I6").Value = 30487.385
$ws.Range("J6").Value = 3000
$ws.Range("K6").Value = 30487.385
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -30374.385
$ws.Range("N6").Value = -3226

$ws.Range("H7").Value = 91.125
$ws.Range("I7").Value = 39.77778
$ws.Range("J7").Value = 157.14285
$ws.Range("K7").Value = 39.77778
$ws.Range("L7").Value = 157.14285
$ws.Range("M7").Value = 73.22221999999999
$ws.Range("N7").Value = -383.14285

$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 4333.3335
$ws.Range("K11").Value = 0
$ws.Range("M11").Value = 4333.3335
$ws.Range("N11").Value = -4613.3335
$ws.Range("L11").ClearContents()

$ws.Range("H12").Value = 13333.333
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 13333.333
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 13333.333
$ws.Range("N12").Value = -13673.333
$ws.Range("M12").ClearContents()

$ws.Range("H13").Value = 18083.334
$ws.Range("J13").Value = 18083.334
$ws.Range("L13").Value = 18083.334
$ws.Range("N13").Value = -18361.334

$ws.Range("H31").Value = 2869.24
$ws.Range("I31").Value = 1733.8334
$ws.Range("J31").Value = 3227.7896
$ws.Range("K31").Value = 1733.8334
$ws.Range("L31").Value = 3227.7896
$ws.Range("M31").Value = -1438.8334
$ws.Range("N31").Value = -3817.7896

$ws.Range("H33").Value = 33639.5
$ws.Range("J33").Value = 33639.5
$ws.Range("L33").Value = 33639.5
$ws.Range("N33").Value = -34397.5

$ws.Range("H34").Value = 2869.24
$ws.Range("I34").Value = 1733.8334
$ws.Range("J34").Value = 3227.7896
$ws.Range("K34").Value = 1733.8334
$ws.Range("L34").Value = 3227.7896
$ws.Range("M34").Value = -1531.8334
$ws.Range("N34").Value = -3631.7896

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 3920
$ws.Range("J22").Value = 3920
$ws.Range("L22").Value = 11760
$ws.Range("N22").Value = -12098

$ws.Range("H27").Value = 3920
$ws.Range("J27").Value = 3920
$ws.Range("L27").Value = 11760
$ws.Range("N27").Value = -11964

$ws.Range("H29").Value = 500
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H34").Value = 2143.3333
$ws.Range("J34").Value = 2452
$ws.Range("L34").Value = 7356
$ws.Range("N34").Value = -7524

$ws.Range("H46").Value = 752788.1
$ws.Range("I46").Value = 1001.5
$ws.Range("J46").Value = 836319.9399999999
$ws.Range("K46").Value = 3004.5
$ws.Range("L46").Value = 2508959.82
$ws.Range("M46").Value = -2913.5
$ws.Range("N46").Value = -2509141.82

$ws.Range("H86").Value = 475
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 475
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = 1425
$ws.Range("N86").Value = -3797
$ws.Range("L86").ClearContents()

$ws.Range("H89").Value = 475
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 475
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = 4275
$ws.Range("N89").Value = -16131
$ws.Range("L89").ClearContents()

$ws.Range("H109").Value = 2109.0667
$ws.Range("I109").Value = 951.1429000000001
$ws.Range("J109").Value = 3122.25
$ws.Range("K109").Value = 2853.4287
$ws.Range("L109").Value = 9366.75
$ws.Range("M109").Value = -1813.4287
$ws.Range("N109").Value = -11446.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5503334
$ws.Range("I11").Value = 5503334
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 5503334
$ws.Range("L11").Value = 0
$ws.Range("N11").Value = -5503195
$ws.Range("M11").ClearContents()

$ws.Range("H12").Value = 23999.6
$ws.Range("J12").Value = 28749.5
$ws.Range("L12").Value = 28749.5
$ws.Range("N12").Value = -29029.5

$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("L13").ClearContents()
$ws.Range("M13").ClearContents()

$ws.Range("H126").Value = 66667984
$ws.Range("I126").Value = 125001176
$ws.Range("J126").Value = 1483.5714
$ws.Range("K126").Value = 375003528
$ws.Range("L126").Value = 4450.7142
$ws.Range("M126").Value = -375001058
$ws.Range("N126").Value = -9390.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1193.3334
$ws.Range("I16").Value = 927.2727
$ws.Range("J16").Value = 1925
$ws.Range("K16").Value = 927.2727
$ws.Range("L16").Value = 1925
$ws.Range("M16").Value = -757.2727
$ws.Range("N16").Value = -2265

$ws.Range("H22").Value = 874
$ws.Range("I22").Value = 690
$ws.Range("J22").Value = 920
$ws.Range("K22").Value = 690
$ws.Range("L22").Value = 920
$ws.Range("M22").Value = -395
$ws.Range("N22").Value = -1510

$ws.Range("H27").Value = 874
$ws.Range("I27").Value = 690
$ws.Range("J27").Value = 920
$ws.Range("K27").Value = 690
$ws.Range("L27").Value = 920
$ws.Range("M27").Value = -583
$ws.Range("N27").Value = -1134

$ws.Range("H46").Value = 896.5789
$ws.Range("I46").Value = 1626.6
$ws.Range("J46").Value = 635.8570999999999
$ws.Range("K46").Value = 1626.6
$ws.Range("L46").Value = 635.8570999999999
$ws.Range("M46").Value = -1438.6
$ws.Range("N46").Value = -1011.8571

$ws.Range("H132").Value = 4392
$ws.Range("I132").Value = 4120.2607
$ws.Range("J132").Value = 5433.6665
$ws.Range("K132").Value = 12360.7821
$ws.Range("L132").Value = 16300.9995
$ws.Range("M132").Value = -9830.7821
$ws.Range("N132").Value = -21360.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 3750925.2
$ws.Range("I6").Value = 750
$ws.Range("K6").Value = 750
$ws.Range("M6").Value = -635

$ws.Range("H11").Value = 1000
$ws.Range("I11").Value = 1000
$ws.Range("K11").Value = 1000
$ws.Range("M11").Value = -858

$ws.Range("H12").Value = 87505.25
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 87505.25
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = 87505.25
$ws.Range("N12").Value = -87789.25
$ws.Range("L12").ClearContents()

$ws.Range("H13").Value = 17900
$ws.Range("I13").Value = 2900
$ws.Range("J13").Value = 25400
$ws.Range("K13").Value = 2900
$ws.Range("L13").Value = 25400
$ws.Range("M13").Value = -2760
$ws.Range("N13").Value = -25680

$ws.Range("H136").Value = 7272.4614
$ws.Range("I136").Value = 6533.778
$ws.Range("J136").Value = 7663.5293
$ws.Range("K136").Value = 19601.334
$ws.Range("L136").Value = 22990.5879
$ws.Range("M136").Value = -17051.334
$ws.Range("N136").Value = -28090.5879
